# Apply the "skill full names added" change:
# Insert a new "Skill Description" column (B) that carries the human-readable
# name of the skill, shifting the existing SFIA Level / Keycode / Description
# columns one place to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B; existing B (SFIA Level), C (Keycode), D (Description)
# shift right to C, D, E respectively.
$ws.Columns("B:B").Insert()

# Header for the newly inserted column.
$ws.Cells.Item(1, 2).Value = "Skill Description"

# Map SkillCode -> full Skill Description text.
$skillNames = @{
    "Autonomy"   = "Autonomy";
    "Influence"  = "Influence";
    "Complexity" = "Complexity";
    "Knowledge"  = "Knowledge";
    "PROF"       = "Portfolio, programme and project support";
    "MEAS"       = "Measurement";
    "METL"       = "Methods and tools";
    "CIPM"       = "Organisational change management";
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Value()
    if ($code -ne $null -and $code -ne "") {
        $ws.Cells.Item($r, 2).Value = $skillNames[$code]
    }
}
